$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append new test rows for netselectivity.sas and volatilityskewness.sas
$ws.Range("A137").Value = "netselectivity"
$ws.Range("B137").Value = "Test netselectivity with BM=SPY and Rf=0.05"
$ws.Range("C137").Value = "netselectivity_test"

$ws.Range("A138").Value = "volatilityskewness"
$ws.Range("B138").Value = "Test volatilityskewnsss with option=VOLATILITY"
$ws.Range("C138").Value = "volatilityskewness_test1"

$ws.Range("A139").Value = "volatilityskewness"
$ws.Range("B139").Value = "Test volatilityskewnsss with option=VARIABILITY"
$ws.Range("C139").Value = "volatilityskewness_test2"

# Update selection to match final state
$ws.Range("A139").Select()
